$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 2222566.8
$ws.Range("I92").Value = 2778107.2
$ws.Range("J92").Value = 405
$ws.Range("K92").Value = 2778107.2
$ws.Range("L92").Value = 405
$ws.Range("M92").Value = -2776859.2
$ws.Range("N92").Value = -2901

$ws.Range("H116").Value = 2017.9348
$ws.Range("I116").Value = 1821.6
$ws.Range("J116").Value = 2251.6667
$ws.Range("K116").Value = 1821.6
$ws.Range("L116").Value = 2251.6667
$ws.Range("M116").Value = 1620.4
$ws.Range("N116").Value = -9135.6667

$ws.Range("H137").Value = 1091.4865
$ws.Range("I137").Value = 1012.92
$ws.Range("J137").Value = 1255.1666
$ws.Range("K137").Value = 3038.76
$ws.Range("L137").Value = 3765.4998
$ws.Range("M137").Value = -488.7599999999998
$ws.Range("N137").Value = -8865.4998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 3800
$ws.Range("I21").Value = 2200
$ws.Range("K21").Value = 2200
$ws.Range("M21").Value = -1826

$ws.Range("H61").Value = 7937991.5
$ws.Range("I61").Value = 7937991.5
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 7937991.5
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -7937779.5
$ws.Range("N61").ClearContents()

$ws.Range("H74").Value = 1267.5667
$ws.Range("I74").Value = 1127.3043
$ws.Range("J74").Value = 1728.4286
$ws.Range("K74").Value = 1127.3043
$ws.Range("L74").Value = 1728.4286
$ws.Range("M74").Value = -253.3043
$ws.Range("N74").Value = -3476.4286

$ws.Range("H77").Value = 1267.5667
$ws.Range("I77").Value = 1127.3043
$ws.Range("J77").Value = 1728.4286
$ws.Range("K77").Value = 5636.5215
$ws.Range("L77").Value = 8642.143
$ws.Range("M77").Value = -1268.5215
$ws.Range("N77").Value = -17378.143

$ws.Range("H132").Value = 920257
$ws.Range("I132").Value = 1167.3269
$ws.Range("J132").Value = 4902979
$ws.Range("K132").Value = 3501.9807
$ws.Range("L132").Value = 14708937
$ws.Range("M132").Value = -971.9807000000001
$ws.Range("N132").Value = -14713997

$ws.Range("H136").Value = 7937991.5
$ws.Range("I136").Value = 7937991.5
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 23813974.5
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -23811424.5
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1225564.8
$ws.Range("I86").Value = 1470.1538
$ws.Range("J86").Value = 3877769.8
$ws.Range("K86").Value = 1470.1538
$ws.Range("L86").Value = 3877769.8
$ws.Range("M86").Value = -347.1538
$ws.Range("N86").Value = -3880015.8

$ws.Range("H89").Value = 1225564.8
$ws.Range("I89").Value = 1470.1538
$ws.Range("J89").Value = 3877769.8
$ws.Range("K89").Value = 7350.769
$ws.Range("L89").Value = 19388849
$ws.Range("M89").Value = -1734.769
$ws.Range("N89").Value = -19400081

$ws.Range("H134").Value = 2383.379
$ws.Range("I134").Value = 845.3333
$ws.Range("J134").Value = 5074.9585
$ws.Range("K134").Value = 2535.9999
$ws.Range("L134").Value = 15224.8755
$ws.Range("M134").Value = -0.9998999999997977
$ws.Range("N134").Value = -20294.8755

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1465471.4
$ws.Range("I31").Value = 1504997.6
$ws.Range("J31").Value = 3000
$ws.Range("K31").Value = 1504997.6
$ws.Range("L31").Value = 3000
$ws.Range("M31").Value = -1504702.6
$ws.Range("N31").Value = -3590

$ws.Range("H34").Value = 1465471.4
$ws.Range("I34").Value = 1504997.6
$ws.Range("J34").Value = 3000
$ws.Range("K34").Value = 1504997.6
$ws.Range("L34").Value = 3000
$ws.Range("M34").Value = -1504795.6
$ws.Range("N34").Value = -3404

$ws.Range("H58").Value = 29412486
$ws.Range("I58").Value = 40000670
$ws.Range("J58").Value = 857.1111
$ws.Range("K58").Value = 40000670
$ws.Range("L58").Value = 857.1111
$ws.Range("M58").Value = -40000467
$ws.Range("N58").Value = -1263.1111

$ws.Range("H86").Value = 43492220
$ws.Range("I86").Value = 71430970
$ws.Range("J86").Value = 31950.777
$ws.Range("K86").Value = 71430970
$ws.Range("L86").Value = 31950.777
$ws.Range("M86").Value = -71429847
$ws.Range("N86").Value = -34196.777

$ws.Range("H89").Value = 43492220
$ws.Range("I89").Value = 71430970
$ws.Range("J89").Value = 31950.777
$ws.Range("K89").Value = 357154850
$ws.Range("L89").Value = 159753.885
$ws.Range("M89").Value = -357149234
$ws.Range("N89").Value = -170985.885

$ws.Range("H134").Value = 1155.5834
$ws.Range("I134").Value = 1024.5714
$ws.Range("J134").Value = 1614.125
$ws.Range("K134").Value = 3073.7142
$ws.Range("L134").Value = 4842.375
$ws.Range("M134").Value = -538.7142000000003
$ws.Range("N134").Value = -9912.375

$ws.Range("H136").Value = 29412486
$ws.Range("I136").Value = 40000670
$ws.Range("J136").Value = 857.1111
$ws.Range("K136").Value = 120002010
$ws.Range("L136").Value = 2571.3333
$ws.Range("M136").Value = -119999460
$ws.Range("N136").Value = -7671.3333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H106").Value = 980
$ws.Range("I106").Value = 980
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 2940
$ws.Range("L106").Value = 0
$ws.Range("M106").Value = -1994
$ws.Range("N106").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4453.553
$ws.Range("I132").Value = 2156.4062
$ws.Range("K132").Value = 6469.2186
$ws.Range("M132").Value = -3939.2186

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 51677308
$ws.Range("I132").Value = 108030296
$ws.Range("J132").Value = 20400.334
$ws.Range("K132").Value = 324090888
$ws.Range("L132").Value = 61201.00199999999
$ws.Range("M132").Value = -324088358
$ws.Range("N132").Value = -66261.00199999999

$ws.Range("H136").Value = 63026772
$ws.Range("I136").Value = 57144656
$ws.Range("J136").Value = 71429790
$ws.Range("K136").Value = 171433968
$ws.Range("L136").Value = 214289370
$ws.Range("M136").Value = -171431418
$ws.Range("N136").Value = -214294470

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()

$ws.Range("H132").Value = 44950.383
$ws.Range("I132").Value = 94573.73
$ws.Range("J132").Value = 8559.933999999999
$ws.Range("K132").Value = 283721.19
$ws.Range("L132").Value = 25679.802
$ws.Range("M132").Value = -281191.19
$ws.Range("N132").Value = -30739.802

$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

$ws.Range("H136").Value = 19233458
$ws.Range("I136").Value = 35716004
$ws.Range("J136").Value = 3823.3333
$ws.Range("K136").Value = 107148012
$ws.Range("L136").Value = 11469.9999
$ws.Range("M136").Value = -107145462
$ws.Range("N136").Value = -16569.9999
